# Update countries & provincias Spain
# - Swap the order of "Catar" and "Nepal" so Nepal's row now comes first
#   (row 38) followed by Catar (row 39), each carrying its own refreshed
#   statistics.
# - Refresh the COVID-19 statistics (Casos totales, Nuevos casos,
#   Casos activos, Recuperados, Muertes hoy, Muertes) for several
#   countries: Estados Unidos, Iran, Nepal, Catar, El Salvador,
#   Sri Lanka, Malta and Vietnam.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap Catar / Nepal order, row 38 becomes Nepal, row 39 becomes Catar ---
$ws.Range("A38").Value = "Nepal"
$ws.Range("A39").Value = "Catar"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 8343244
$ws.Range("C4").Value = 579
$ws.Range("D4").Value = 5432457
$ws.Range("E4").Value = 2686503
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 224284

# --- Row 16: Iran ---
$ws.Range("B16").Value = 530380
$ws.Range("C16").Value = 3890
$ws.Range("D16").Value = 427400
$ws.Range("E16").Value = 72605
$ws.Range("G16").Value = 252
$ws.Range("H16").Value = 30375

# --- Row 38: Nepal (new data) ---
$ws.Range("B38").Value = 132246
$ws.Range("C38").Value = 2942
$ws.Range("D38").Value = 92166
$ws.Range("E38").Value = 39341
$ws.Range("G38").Value = 12
$ws.Range("H38").Value = 739

# --- Row 39: Catar (carries over previous Catar data) ---
$ws.Range("B39").Value = 129431
$ws.Range("C39").Value = 204
$ws.Range("D39").Value = 126406
$ws.Range("E39").Value = 2801
$ws.Range("G39").Value = 1
$ws.Range("H39").Value = 224

# --- Row 83: El Salvador ---
$ws.Range("B83").Value = 31666
$ws.Range("C83").Value = 210
$ws.Range("D83").Value = 27000
$ws.Range("E83").Value = 3744

# --- Row 126: Sri Lanka ---
$ws.Range("B126").Value = 5497
$ws.Range("C126").Value = 22
$ws.Range("E126").Value = 2081

# --- Row 139: Malta ---
$ws.Range("B139").Value = 4628
$ws.Range("C139").Value = 142
$ws.Range("D139").Value = 3236
$ws.Range("E139").Value = 1347

# --- Row 168: Vietnam ---
$ws.Range("B168").Value = 1134
$ws.Range("C168").Value = 8
$ws.Range("E168").Value = 68
